$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Row, $StartCol, $EndCol, $Values)
    $arr = New-Object "object[,]" 1, $Values.Count
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $ws.Range($ws.Cells.Item($Row, $StartCol), $ws.Cells.Item($Row, $EndCol)).Value2 = $arr
}

# NOTE: this engine requires POSITIONAL params: Set-RowValues <Row> <StartCol> <EndCol> <Values>

# --- Re-order same-date matches (swap/rotate F:V content across row groups) ---
# Group [20, 21]
Set-RowValues 20 6 22 @("Vitoria", 2, "Londrina", 0, 2.13, "23/04/2023 23:12", 1.84, "28/04/2023 23:59", 3.31, "23/04/2023 23:12", 3.45, "28/04/2023 23:59", 3.45, "23/04/2023 23:12", 4.94, "28/04/2023 23:59", "https://www.betexplorer.com/football/brazil/serie-b/vitoria-londrina/W80E7CgF/")
Set-RowValues 21 6 22 @("Criciuma", 1, "Avai", 0, 1.73, "22/04/2023 21:12", 1.76, "28/04/2023 23:49", 3.37, "22/04/2023 21:12", 3.49, "28/04/2023 23:59", 5.21, "22/04/2023 21:12", 5.46, "28/04/2023 23:59", "https://www.betexplorer.com/football/brazil/serie-b/criciuma-avai/xd4Opkn2/")

# Group [31, 32]
Set-RowValues 31 6 22 @("Ponte Preta", 2, "Botafogo SP", 0, 2.15, "29/04/2023 22:13", 2.03, "02/05/2023 23:23", 3.06, "29/04/2023 22:13", 3.25, "02/05/2023 23:23", 3.72, "29/04/2023 22:13", 4.28, "02/05/2023 23:23", "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-botafogo-sp/xAqcal1e/")
Set-RowValues 32 6 22 @("CRB", 0, "Sampaio Correa", 0, 1.92, "29/04/2023 22:13", 1.86, "02/05/2023 23:59", 3.33, "29/04/2023 22:13", 3.51, "02/05/2023 23:57", 4.47, "29/04/2023 22:13", 4.68, "02/05/2023 23:59", "https://www.betexplorer.com/football/brazil/serie-b/crb-sampaio-correa/pf7UaLi6/")

# Group [40, 41, 42]
Set-RowValues 40 6 22 @("Tombense", 0, "Avai", 1, 2.2, "04/05/2023 02:42", 2.03, "06/05/2023 23:13", 3.21, "04/05/2023 02:42", 3.32, "06/05/2023 23:13", 3.4, "04/05/2023 02:42", 4.16, "06/05/2023 23:13", "https://www.betexplorer.com/football/brazil/serie-b/tombense-avai/xlW7DvM5/")
Set-RowValues 41 6 22 @("Sampaio Correa", 2, "Juventude", 1, 2.3, "04/05/2023 00:12", 2, "06/05/2023 23:13", 3.08, "04/05/2023 00:12", 3.44, "06/05/2023 23:13", 3.31, "04/05/2023 00:12", 4.11, "06/05/2023 23:13", "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-esporte-clube-juventude/UNSaF0ig/")
Set-RowValues 42 6 22 @("Mirassol", 1, "Vila Nova FC", 1, 2.07, "06/05/2023 15:09", 1.95, "06/05/2023 23:09", 3.06, "06/05/2023 15:09", 3.3, "06/05/2023 23:10", 3.96, "06/05/2023 15:09", 4.55, "06/05/2023 23:14", "https://www.betexplorer.com/football/brazil/serie-b/mirassol-vila-nova-fc/YeZjHMMt/")

# Group [53, 54]
Set-RowValues 53 6 22 @("Guarani", 2, "Sampaio Correa", 0, 1.72, "10/05/2023 02:41", 1.88, "13/05/2023 21:58", 3.87, "10/05/2023 02:41", 3.38, "13/05/2023 21:58", 5.1, "10/05/2023 02:41", 4.78, "13/05/2023 21:58", "https://www.betexplorer.com/football/brazil/serie-b/guarani-sampaio-correa/b3KGBIjI/")
Set-RowValues 54 6 22 @("Novorizontino", 0, "CRB", 1, 1.69, "10/05/2023 02:42", 1.76, "13/05/2023 21:55", 3.66, "10/05/2023 02:42", 3.56, "13/05/2023 21:12", 4.94, "10/05/2023 02:42", 5.33, "13/05/2023 21:55", "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-crb/COkWD8Xo/")

# Group [57, 58]
Set-RowValues 57 6 22 @("ABC", 1, "Botafogo SP", 2, 2.18, "10/05/2023 02:42", 2.17, "14/05/2023 22:57", 3.11, "10/05/2023 02:42", 3.06, "14/05/2023 22:54", 3.81, "10/05/2023 02:42", 4.08, "14/05/2023 22:57", "https://www.betexplorer.com/football/brazil/serie-b/abc-botafogo-sp/zRti9pnH/")
Set-RowValues 58 6 22 @("Vitoria", 2, "Atletico GO", 3, 2.04, "11/05/2023 00:12", 1.95, "14/05/2023 22:59", 3.28, "11/05/2023 00:12", 3.38, "14/05/2023 22:49", 3.76, "11/05/2023 00:12", 4.39, "14/05/2023 22:59", "https://www.betexplorer.com/football/brazil/serie-b/vitoria-atletico-go/6okzDSnh/")

# Group [64, 65]
Set-RowValues 64 6 22 @("Ituano", 0, "Novorizontino", 2, 2.65, "15/05/2023 01:42", 2.71, "20/05/2023 21:59", 3.13, "15/05/2023 01:42", 3.08, "20/05/2023 21:59", 2.75, "15/05/2023 01:42", 2.94, "20/05/2023 21:50", "https://www.betexplorer.com/football/brazil/serie-b/ituano-novorizontino/KlF9RR1o/")
Set-RowValues 65 6 22 @("Sampaio Correa", 1, "ABC", 0, 1.93, "14/05/2023 23:12", 1.91, "20/05/2023 21:52", 3.34, "14/05/2023 23:12", 3.47, "20/05/2023 21:57", 4.4, "14/05/2023 23:12", 4.47, "20/05/2023 21:57", "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-abc/f77MOPo4/")

# Group [85, 86]
Set-RowValues 85 6 22 @("Mirassol", 1, "Criciuma", 0, 2.04, "26/05/2023 03:42", 2, "28/05/2023 20:21", 3.36, "26/05/2023 03:42", 3.31, "28/05/2023 20:21", 4.05, "26/05/2023 03:42", 4.3, "28/05/2023 20:21", "https://www.betexplorer.com/football/brazil/serie-b/mirassol-criciuma/td9WxG6E/")
Set-RowValues 86 6 22 @("Ceara", 0, "Novorizontino", 3, 2.05, "26/05/2023 03:42", 2.09, "28/05/2023 20:20", 3.22, "26/05/2023 03:42", 3.18, "28/05/2023 20:29", 3.8, "26/05/2023 03:42", 4.15, "28/05/2023 20:29", "https://www.betexplorer.com/football/brazil/serie-b/ceara-novorizontino/2BcuyfyR/")

# Group [90, 91]
Set-RowValues 90 6 22 @("Botafogo SP", 0, "Tombense", 0, 2.09, "31/05/2023 05:12", 2.05, "02/06/2023 23:58", 3.28, "31/05/2023 05:12", 3.36, "02/06/2023 23:58", 3.61, "31/05/2023 05:12", 4.03, "02/06/2023 23:59", "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-tombense/IsKAjZzE/")
Set-RowValues 91 6 22 @("Ceara", 2, "Chapecoense-SC", 0, 1.7, "31/05/2023 05:12", 1.61, "02/06/2023 23:34", 3.48, "31/05/2023 05:12", 3.87, "02/06/2023 23:34", 5.2, "31/05/2023 05:12", 6.37, "02/06/2023 23:34", "https://www.betexplorer.com/football/brazil/serie-b/ceara-chapecoense-sc/Qm9ppXDs/")

# Group [92, 93]
Set-RowValues 92 6 22 @("Vitoria", 3, "Ituano", 0, 1.69, "31/05/2023 05:12", 1.78, "03/06/2023 02:24", 3.53, "31/05/2023 05:12", 3.51, "03/06/2023 02:24", 5.21, "31/05/2023 05:12", 5.29, "03/06/2023 02:24", "https://www.betexplorer.com/football/brazil/serie-b/vitoria-ituano/jeIIlD5Q/")
Set-RowValues 93 6 22 @("Criciuma", 3, "Atletico GO", 0, 2.06, "31/05/2023 05:12", 2.11, "03/06/2023 02:13", 3.24, "31/05/2023 05:12", 3.1, "03/06/2023 02:13", 3.73, "31/05/2023 05:12", 4.21, "03/06/2023 02:13", "https://www.betexplorer.com/football/brazil/serie-b/criciuma-atletico-go/tQzSRkbD/")

# Group [112, 113]
Set-RowValues 112 6 22 @("Novorizontino", 1, "Sampaio Correa", 0, 1.56, "08/06/2023 14:42", 1.61, "10/06/2023 21:22", 3.64, "08/06/2023 14:42", 3.81, "10/06/2023 21:22", 6.5, "08/06/2023 14:42", 6.55, "10/06/2023 21:22", "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-sampaio-correa/QNPgB9Ct/")
Set-RowValues 113 6 22 @("Ituano", 1, "Atletico GO", 1, 2.51, "07/06/2023 02:42", 2.28, "10/06/2023 21:51", 3.03, "07/06/2023 02:42", 3.22, "10/06/2023 21:59", 3.02, "07/06/2023 02:42", 3.5, "10/06/2023 21:59", "https://www.betexplorer.com/football/brazil/serie-b/ituano-atletico-go/jVYRFk4P/")

# Group [176, 177]
Set-RowValues 176 6 22 @("Avai", 2, "Sampaio Correa", 0, 2.2, "15/07/2023 02:42", 2.05, "20/07/2023 02:29", 3.05, "15/07/2023 02:42", 2.98, "20/07/2023 02:29", 3.82, "15/07/2023 02:42", 4.71, "20/07/2023 02:29", "https://www.betexplorer.com/football/brazil/serie-b/avai-sampaio-correa/fPQDLTgU/")
Set-RowValues 177 6 22 @("ABC", 0, "Guarani", 1, 2.6, "15/07/2023 16:12", 2.6, "20/07/2023 02:29", 2.85, "15/07/2023 16:12", 2.81, "20/07/2023 02:26", 3.08, "15/07/2023 16:12", 3.42, "20/07/2023 02:29", "https://www.betexplorer.com/football/brazil/serie-b/abc-guarani/KS9yy9gH/")

# Group [276, 277]
Set-RowValues 276 6 22 @("Sampaio Correa", 2, "Chapecoense-SC", 0, 2.08, "11/09/2023 02:42", 2.12, "16/09/2023 21:50", 2.96, "11/09/2023 02:42", 2.89, "16/09/2023 21:50", 4.43, "11/09/2023 02:42", 4.62, "16/09/2023 21:50", "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-chapecoense-sc/6gaZd5YQ/")
Set-RowValues 277 6 22 @("Botafogo SP", 1, "Atletico GO", 0, 2.92, "11/09/2023 08:12", 3.52, "16/09/2023 21:52", 2.79, "11/09/2023 08:12", 2.88, "16/09/2023 21:50", 2.92, "11/09/2023 08:12", 2.48, "16/09/2023 21:50", "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-atletico-go/M5cNa7l8/")

# Group [286, 287]
Set-RowValues 286 6 22 @("CRB", 1, "Guarani", 0, 2.09, "18/09/2023 20:13", 2.14, "23/09/2023 21:53", 3.03, "18/09/2023 20:13", 3.06, "23/09/2023 21:53", 4.27, "18/09/2023 20:13", 4.15, "23/09/2023 21:53", "https://www.betexplorer.com/football/brazil/serie-b/crb-guarani/KMGp3te0/")
Set-RowValues 287 6 22 @("Chapecoense-SC", 1, "Ceara", 1, 2.55, "19/09/2023 01:12", 2.4, "23/09/2023 21:58", 2.89, "19/09/2023 01:12", 2.91, "23/09/2023 21:59", 3.27, "19/09/2023 01:12", 3.67, "23/09/2023 21:59", "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-ceara/OfNy52Qm/")

# Group [307, 308]
Set-RowValues 307 6 22 @("Botafogo SP", 0, "Avai", 1, 2.02, "02/10/2023 06:12", 2.05, "07/10/2023 21:58", 3.01, "02/10/2023 06:12", 3.04, "07/10/2023 21:58", 4.62, "02/10/2023 06:12", 4.58, "07/10/2023 21:58", "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-avai/hnaoLEGo/")
Set-RowValues 308 6 22 @("Sampaio Correa", 1, "Novorizontino", 1, 3.32, "02/10/2023 06:12", 3.03, "07/10/2023 21:59", 2.94, "02/10/2023 06:12", 2.89, "07/10/2023 21:51", 2.39, "02/10/2023 06:12", 2.8, "07/10/2023 21:59", "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-novorizontino/S4ibIC04/")

# Group [315, 316]
Set-RowValues 315 6 22 @("Novorizontino", 0, "Tombense", 0, 1.51, "08/10/2023 20:12", 1.52, "14/10/2023 21:54", 3.89, "08/10/2023 20:12", 4.03, "14/10/2023 21:54", 6.63, "08/10/2023 20:12", 7.51, "14/10/2023 21:54", "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-tombense/6eCXWBFc/")
Set-RowValues 316 6 22 @("Ituano", 0, "CRB", 0, 2.37, "07/10/2023 21:12", 2.3, "14/10/2023 21:59", 2.95, "07/10/2023 21:12", 2.99, "14/10/2023 21:59", 3.55, "07/10/2023 21:12", 3.79, "14/10/2023 21:59", "https://www.betexplorer.com/football/brazil/serie-b/ituano-crb/S4DTXi0i/")

# Group [322, 323]
Set-RowValues 322 6 22 @("Atletico GO", 3, "ABC", 1, 1.35, "15/10/2023 20:12", 1.39, "19/10/2023 23:51", 4.58, "15/10/2023 20:12", 4.61, "19/10/2023 23:58", 9.11, "15/10/2023 20:12", 9.75, "19/10/2023 23:58", "https://www.betexplorer.com/football/brazil/serie-b/atletico-go-abc/YuRGMRij/")
Set-RowValues 323 6 22 @("Tombense", 2, "Vila Nova FC", 1, 2.81, "14/10/2023 21:13", 3.42, "19/10/2023 23:47", 2.86, "14/10/2023 21:13", 2.91, "19/10/2023 23:47", 2.97, "14/10/2023 21:13", 2.52, "19/10/2023 23:47", "https://www.betexplorer.com/football/brazil/serie-b/tombense-vila-nova-fc/8I9Nu7bS/")

# Group [329, 330]
Set-RowValues 329 6 22 @("Sport Recife", 2, "Chapecoense-SC", 1, 1.58, "17/10/2023 00:12", 1.47, "22/10/2023 22:52", 3.68, "17/10/2023 00:12", 4.18, "22/10/2023 22:59", 6.14, "17/10/2023 00:12", 8.48, "22/10/2023 22:59", "https://www.betexplorer.com/football/brazil/serie-b/sport-recife-chapecoense-sc/M17Br9E9/")
Set-RowValues 330 6 22 @("Avai", 1, "Ceara", 0, 2.47, "16/10/2023 03:12", 2.34, "22/10/2023 22:51", 3.1, "16/10/2023 03:12", 3.13, "22/10/2023 22:51", 3.01, "16/10/2023 03:12", 3.48, "22/10/2023 22:50", "https://www.betexplorer.com/football/brazil/serie-b/avai-ceara/CWPKLo7d/")

# Group [334, 335]
Set-RowValues 334 6 22 @("Guarani", 0, "Botafogo SP", 1, 1.64, "21/10/2023 22:12", 1.7, "28/10/2023 02:05", 3.41, "21/10/2023 22:12", 3.43, "28/10/2023 02:17", 6.14, "21/10/2023 22:12", 6.36, "28/10/2023 02:24", "https://www.betexplorer.com/football/brazil/serie-b/guarani-botafogo-sp/QZkfzoMk/")
Set-RowValues 335 6 22 @("Novorizontino", 2, "Ponte Preta", 0, 1.53, "24/10/2023 01:12", 1.39, "28/10/2023 02:21", 3.81, "24/10/2023 01:12", 4.57, "28/10/2023 02:28", 6.57, "24/10/2023 01:12", 10.13, "28/10/2023 02:28", "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-ponte-preta/EeWldSr9/")

# Group [336, 337]
Set-RowValues 336 6 22 @("Criciuma", 3, "Sampaio Correa", 0, 1.7, "21/10/2023 23:42", 1.58, "28/10/2023 20:26", 3.42, "21/10/2023 23:42", 3.72, "28/10/2023 20:26", 5.95, "21/10/2023 23:42", 7.22, "28/10/2023 20:26", "https://www.betexplorer.com/football/brazil/serie-b/criciuma-sampaio-correa/bml2YOj2/")
Set-RowValues 337 6 22 @("Ituano", 0, "Mirassol", 0, 3.08, "24/10/2023 01:12", 3.71, "28/10/2023 20:24", 2.89, "24/10/2023 01:12", 2.9, "28/10/2023 20:24", 2.57, "24/10/2023 01:12", 2.39, "28/10/2023 20:24", "https://www.betexplorer.com/football/brazil/serie-b/ituano-mirassol/z5Xpc8T2/")

# --- Append new rows for additional matches (355-357) ---
# Row 355
Set-RowValues 355 1 22 @(354, "brazil", "serie-b", "2023", 45241.875, "Chapecoense-SC", 2, "Botafogo SP", 0, 2.01, "07/11/2023 23:12", 1.85, "11/11/2023 20:56", 3.1, "07/11/2023 23:12", 3.32, "11/11/2023 20:56", 4.44, "07/11/2023 23:12", 5.14, "11/11/2023 20:56", "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-botafogo-sp/jwktn3xa/")
$ws.Cells.Item(355, 1).Style = $ws.Cells.Item(354, 1).Style
$ws.Cells.Item(355, 5).Style = $ws.Cells.Item(354, 5).Style

# Row 356
Set-RowValues 356 1 22 @(355, "brazil", "serie-b", "2023", 45241.875, "Avai", 1, "CRB", 1, 2.33, "07/11/2023 23:12", 2.15, "11/11/2023 20:55", 3.04, "07/11/2023 23:12", 3.23, "11/11/2023 20:31", 3.31, "07/11/2023 23:12", 3.84, "11/11/2023 20:55", "https://www.betexplorer.com/football/brazil/serie-b/avai-crb/dKHSpgtP/")
$ws.Cells.Item(356, 1).Style = $ws.Cells.Item(355, 1).Style
$ws.Cells.Item(356, 5).Style = $ws.Cells.Item(355, 5).Style

# Row 357
Set-RowValues 357 1 22 @(356, "brazil", "serie-b", "2023", 45241.91666666666, "Tombense", 0, "Ponte Preta", 1, 1.93, "04/11/2023 23:43", 1.67, "11/11/2023 21:51", 3.22, "04/11/2023 23:43", 3.61, "11/11/2023 21:51", 4.66, "04/11/2023 23:43", 6.18, "11/11/2023 21:51", "https://www.betexplorer.com/football/brazil/serie-b/tombense-ponte-preta/UTmpoNi5/")
$ws.Cells.Item(357, 1).Style = $ws.Cells.Item(356, 1).Style
$ws.Cells.Item(357, 5).Style = $ws.Cells.Item(356, 5).Style
